$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (the "R40" rule row) has its Rule-id cell (B11) changed from the
# text "R40" to the text "1". Excel would normally auto-convert a bare
# "1" typed into a General-formatted cell into a number, so we briefly
# stash the cell's existing formatting, force a Text number format while
# we assign the new value (so it is stored as a genuine text/string
# cell), then restore the original formatting and clean up the scratch
# cell used to hold it.
$target = $ws.Range("B11")
$scratch = $ws.Range("ZZ1000")

$target.Copy($scratch)          # remember B11's current formatting
$target.NumberFormat = "@"      # force text storage
$target.Value = "1"             # new rule id, stored as text "1"
$scratch.Copy()
$target.PasteSpecial(-4122)     # xlPasteFormats: restore original formatting
$scratch.Delete()               # remove the scratch cell entirely
$excel.CutCopyMode = $false
